$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Índice"
$ws.Range("B1").Value = "Distancia"
$ws.Range("C1").Value = "max"
$ws.Range("D1").Value = "min"
$ws.Range("E1").Value = "Tempo"

$data = @(
    @(0, 2242.833333333333, 2376, 2098, 0.03446881771087647),
    @(1, 2247.266666666667, 2331, 2165, 0.03574072519938151),
    @(2, 1987.2, 2254, 1784, 0.0419516642888387),
    @(3, 2313.8, 2447, 2199, 0.04031805197397868),
    @(4, 1974.966666666667, 2133, 1768, 0.03744359811147054),
    @(5, 2146.7, 2266, 2051, 0.03668584823608399),
    @(6, 2228.233333333333, 2355, 2136, 0.03848103682200114),
    @(7, 2261.033333333333, 2476, 2100, 0.04106341203053792),
    @(8, 2425.133333333333, 2588, 2267, 0.03517295519510905),
    @(9, 1968.733333333333, 2103, 1938, 0.0356726884841919)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r++
}
